$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dates that Excel would otherwise mis-parse as real dates (day <= 12) need a
# leading apostrophe so they are stored as literal text, matching the original
# inline-string cells. We restore the style afterwards so no stray number
# format / quote-prefix formatting is left on the cell.
function Set-LiteralDate($addr, $text) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).Style = "Normal"
}

# Row 3
Set-LiteralDate "A3" "28-07-2022"
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Row 4
Set-LiteralDate "A4" "01-08-2022"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

# Row 5
Set-LiteralDate "A5" "04-08-2022"

# Row 6
Set-LiteralDate "A6" "08-08-2022"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("H6").Value = 0

# Row 7
Set-LiteralDate "A7" "11-08-2022"

# Row 8
Set-LiteralDate "A8" "15-08-2022"

# Row 9
Set-LiteralDate "A9" "18-08-2022"

# Row 10
Set-LiteralDate "A10" "22-08-2022"

# Row 11
Set-LiteralDate "A11" "25-08-2022"
$ws.Range("D11").Value = 1
$ws.Range("G11").Value = 1

# Row 12
Set-LiteralDate "A12" "29-08-2022"
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("H12").Value = 0

# Row 13
Set-LiteralDate "A13" "01-09-2022"
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("H13").Value = 0

# Row 14
Set-LiteralDate "A14" "05-09-2022"

# Row 15
Set-LiteralDate "A15" "08-09-2022"

# Row 16
Set-LiteralDate "A16" "12-09-2022"

# Row 17
Set-LiteralDate "A17" "15-09-2022"

# Row 18
Set-LiteralDate "A18" "19-09-2022"

# Row 19
Set-LiteralDate "A19" "22-09-2022"

# Row 20
Set-LiteralDate "A20" "26-09-2022"

# Row 21
Set-LiteralDate "A21" "29-09-2022"
